# Swap the presentation's theme color palette from the "Integral" / Red
# Violet scheme over to the standard Office color scheme (dk1/lt1/dk2/lt2/
# accent1-6/hlink/folHlink), matching the 12-slot MsoThemeColorSchemeIndex
# ordering exposed through ThemeColorScheme.
#
# RGB() is not available in this host, so pack r,g,b into the same
# 0xBBGGRR long that VBA's RGB() produces.
function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = (RGBVal 0x00 0x00 0x00)   # dk1
$tcs.Item(2).RGB  = (RGBVal 0xFF 0xFF 0xFF)   # lt1
$tcs.Item(3).RGB  = (RGBVal 0x44 0x54 0x6A)   # dk2
$tcs.Item(4).RGB  = (RGBVal 0xE7 0xE6 0xE6)   # lt2
$tcs.Item(5).RGB  = (RGBVal 0x5B 0x9B 0xD5)   # accent1
$tcs.Item(6).RGB  = (RGBVal 0xED 0x7D 0x31)   # accent2
$tcs.Item(7).RGB  = (RGBVal 0xA5 0xA5 0xA5)   # accent3
$tcs.Item(8).RGB  = (RGBVal 0xFF 0xC0 0x00)   # accent4
$tcs.Item(9).RGB  = (RGBVal 0x44 0x72 0xC4)   # accent5
$tcs.Item(10).RGB = (RGBVal 0x70 0xAD 0x47)   # accent6
$tcs.Item(11).RGB = (RGBVal 0x05 0x63 0xC1)   # hlink
$tcs.Item(12).RGB = (RGBVal 0x95 0x4F 0x72)   # folHlink
